$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the dataset. It belongs logically
# before the current row 348, so insert a new row there; this pushes the
# existing rows 348:364 down to 349:365 (shifting all of their data intact).
$ws.Rows("348:348").Insert()

# Populate the newly inserted row 348 with the new record's data.
$ws.Range("A348").Value = 7
$ws.Range("B348").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C348").Value = "Ñuble"
$ws.Range("D348").Value = 45147
$ws.Range("E348").Value = 16
$ws.Range("F348").Value = 100112043
$ws.Range("G348").Value = "Pepino ensalada"
$ws.Range("H348").Value = "Sin especificar"
$ws.Range("I348").Value = "Primera"
$ws.Range("J348").Value = 80
$ws.Range("K348").Value = 12000
$ws.Range("L348").Value = 12000
$ws.Range("M348").Value = 12000
$ws.Range("N348").Value = "`$/caja 60 unidades"
$ws.Range("O348").Value = "Región de Arica y Parinacota"
$ws.Range("P348").Value = 200
$ws.Range("Q348").Value = 60
$ws.Range("R348").Value = "Hortaliza"
